# Update the threshold values on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.6
$ws.Range("B3").Value = 5.7
$ws.Range("B4").Value = 0.9
$ws.Range("C4").Value = 1.2

# Leave the selection where the user last clicked (C4), matching the
# saved sheet view state.
$ws.Range("C4").Select() | Out-Null
